# CotSoc added in amounts. Aggregates widget fixed.
#
# - On "amounts": swap the cotsal_noncontrib / cotpat_noncontrib labels on
#   rows 4 and 7, and append two new aggregate rows (cotsoc_noncontrib,
#   cotsoc_contrib).
# - On "montants": swap the same labels on rows 4 and 7, then fill in the
#   previously-blank "aggregates" rows 9 and 10 with the cotsoc_noncontrib /
#   cotsoc_contrib formulas (and a warning note in J9).
# - Update the active sheet / selections to match: "amounts" becomes the
#   selected tab (cell B9 selected), "montants" keeps A9:I10 selected.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "amounts"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("amounts")

# Rows 4 and 7 had their labels swapped (cotsal_noncontrib <-> cotpat_noncontrib)
$ws1.Range("A4").Value = "cotpat_noncontrib"
$ws1.Range("A7").Value = "cotsal_noncontrib"

# New rows 8 and 9: copy formatting from existing rows so styles/number
# formats match, then overwrite with the new labels/values.
$ws1.Range("A4:I4").Copy($ws1.Range("A8:I8"))
$ws1.Range("A3:I3").Copy($ws1.Range("A9:I9"))

$ws1.Range("A8").Value = "cotsoc_noncontrib"
$ws1.Range("B8").Value = 81241000000
$ws1.Range("C8").Value = 83504000000
$ws1.Range("D8").Value = 86067000000
$ws1.Range("E8").Value = 89023000000
$ws1.Range("F8").Value = 91884000000
$ws1.Range("G8").Value = 94363000000
$ws1.Range("H8").Value = 94387000000
$ws1.Range("I8").Value = 96944000000

$ws1.Range("A9").Value = "cotsoc_contrib"
$ws1.Range("B9").Value = 57801000000
$ws1.Range("C9").Value = 59527000000
$ws1.Range("D9").Value = 61982000000
$ws1.Range("E9").Value = 65130000000
$ws1.Range("F9").Value = 66686000000
$ws1.Range("G9").Value = 67857000000
$ws1.Range("H9").Value = 68049000000
$ws1.Range("I9").Value = 69422000000

# ---------------------------------------------------------------------
# Sheet "montants"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("montants")

# Same label swap as on "amounts"
$ws2.Range("A4").Value = "cotpat_noncontrib"
$ws2.Range("A7").Value = "cotsal_noncontrib"

# Rows 9 and 10 were blank placeholder rows; copy formatting from row 6
# (same label/data style) then fill with the new aggregate formulas.
$ws2.Range("A6:I6").Copy($ws2.Range("A9:I9"))
$ws2.Range("A6:I6").Copy($ws2.Range("A10:I10"))

$ws2.Range("A9").Value = "cotsoc_noncontrib"
$ws2.Range("B9").Formula = "=B4+B7"
$ws2.Range("C9").Formula = "=C4+C7"
$ws2.Range("D9").Formula = "=D4+D7"
$ws2.Range("E9").Formula = "=E4+E7"
$ws2.Range("F9").Formula = "=F4+F7"
$ws2.Range("G9").Formula = "=G4+G7"
$ws2.Range("H9").Formula = "=H4+H7"
$ws2.Range("I9").Formula = "=I4+I7"
$ws2.Range("J9").Value = "Attention: totaux hors non salariés"

$ws2.Range("A10").Value = "cotsoc_contrib"
$ws2.Range("B10").Formula = "=B3+B6"
$ws2.Range("C10").Formula = "=C3+C6"
$ws2.Range("D10").Formula = "=D3+D6"
$ws2.Range("E10").Formula = "=E3+E6"
$ws2.Range("F10").Formula = "=F3+F6"
$ws2.Range("G10").Formula = "=G3+G6"
$ws2.Range("H10").Formula = "=H3+H6"
$ws2.Range("I10").Formula = "=I3+I6"

# ---------------------------------------------------------------------
# Selections / active sheet: select montants first so the final .Select()
# (on amounts) leaves "amounts" as the active/selected tab, matching the
# target workbook view state.
# ---------------------------------------------------------------------
$ws2.Range("A9:I10").Select()
$ws1.Range("B9").Select()
